# Apply the reordering of names in the A1:E6 range on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @(
    @("Loic",    "Antoine", "Geoffroy", "Damien",   "Rıdvan"),
    @("Alisher", "Minh",    "Christian","Mustafa",  "Claire"),
    @("Benjamin","Rachel",  "Ezgi",     "Ness",     "Mehmet"),
    @("Colin",   "Laura",   "Georgina", "Alper",    "Adrien"),
    @("Aadel",   "Hui",     "Rıdvan",   "Mathieu",  "Rıdvan"),
    @("Kyllian", $null,     $null,      $null,      $null)
)

for ($r = 0; $r -lt $values.Length; $r++) {
    $row = $values[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $val = $row[$c]
        $cell = $ws.Cells.Item($r + 1, $c + 1)
        if ($null -eq $val) {
            $cell.Value = ""
        } else {
            $cell.Value = $val
        }
    }
}
